# "fix pull image for synch chapter"
#
# 1) Remove the leftover "pull" callout shapes (a connector + three ovals)
#    from the Synchronization chapter slide (sldId=278) that were left
#    behind from an earlier edit.
# 2) Refresh the master/layout/notes-master "datetimeFigureOut" footer
#    fields (PowerPoint re-stamps these to the current date whenever the
#    deck is saved).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Delete the stray shapes on the slide with SlideID 278.
# ---------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.SlideID -eq 278) {
        $target = $candidate
        break
    }
}

if ($target -ne $null) {
    $namesToRemove = @("Straight Connector 39", "Oval 40", "Oval 57", "Oval 58")
    foreach ($shapeName in $namesToRemove) {
        for ($i = $target.Shapes.Count; $i -ge 1; $i--) {
            $sh = $target.Shapes.Item($i)
            if ($sh.Name -eq $shapeName) {
                $sh.Delete()
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Update the "last saved" date footer text to 2/9/25.
# ---------------------------------------------------------------------
$newDate = "2/9/25"

# Slide master date placeholder.
foreach ($sh in $p.SlideMaster.Shapes) {
    if ($sh.HasTextFrame) {
        if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Every slide layout's date placeholder.
foreach ($cl in $p.SlideMaster.CustomLayouts) {
    foreach ($sh in $cl.Shapes) {
        if ($sh.HasTextFrame) {
            if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Notes master date placeholder (only reachable through HeadersFooters).
$p.NotesMaster.HeadersFooters.DateAndTime.Text = $newDate
